$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at row 79, shifting existing rows 79:140 down to 84:145
$ws.Rows("79:83").Insert()

# Fill in the new rows with the new content
$newRows = @(
    @("Thank you for everything.", "いろいろおせわになりました。"),
    @("Please take care of yourself.", "体に気をつけてください。|からだにきをつけてください。"),
    @("I am looking forward to seeing you.", "お会いできるのを楽しみにしています。|おあいできるをたのしみにしています。"),
    @("Congratulations on...", "～おめでとう（ございます）。"),
    @("Happy Birthday.", "（お）たんじょうびおめでとう。")
)

$r = 79
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}
